$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.654.40"
$ws.Range("E2").Value = "  -3.19%  "

# Row 3
$ws.Range("D3").Value = "1.739.89"
$ws.Range("E3").Value = "  -5.57%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -10.18%  "

# Row 6
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4906"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -7.73%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.27"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -8.43%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2566"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -17.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06065"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -11.98%  "

# Row 11
$ws.Range("D11").Value = "1.741.52"
$ws.Range("E11").Value = "  -5.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06839"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -12.70%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.77"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -20.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.442"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -11.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.74"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -15.47%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.5598"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -26.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.05%  "

# Row 18
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("D19").Value = "25.685.31"
$ws.Range("E19").Value = "  -3.17%  "

# Row 20
$ws.Range("E20").Value = "  -18.67%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006560"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -17.37%  "

# Row 22
$ws.Range("D22").Value = "1.962.17"
$ws.Range("E22").Value = "  -6.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.029"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -12.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.894"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -15.22%  "

# Row 25
$ws.Range("E25").Value = "  -16.87%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.465"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -13.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.818"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -16.79%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.66"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -13.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "101.01"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -9.29%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07964"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.54%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.684"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -14.00%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.391"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -17.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04393"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -9.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9996"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.613"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -11.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9718"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -14.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5898"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -19.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.646"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -14.77%  "

# Row 40
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.07"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01503"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -12.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.839"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -21.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.130"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -12.95%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3727"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -22.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7174"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -20.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05212"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -10.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1084"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -12.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -14.42%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.786"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -24.25%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.67"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -14.28%  "
